$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.731.68"
$ws.Range("E2").Value = "  +0.88%  "
$ws.Range("D3").Value = "'1.773.84"
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'327.16"
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "'0.4607"
$ws.Range("E7").Value = "  +3.57%  "
$ws.Range("D8").Value = "'0.3580"
$ws.Range("E8").Value = "  -0.56%  "
$ws.Range("D9").Value = "'0.07490"
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").Value = "'41.88"
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").Value = "'20.83"
$ws.Range("E13").Value = "  +1.16%  "
$ws.Range("D14").Value = "'6.045"
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("D15").Value = "'7.225"
$ws.Range("E15").Value = "  +1.33%  "
$ws.Range("D16").Value = "'1.770.92"
$ws.Range("E16").Value = "  +1.05%  "
$ws.Range("D17").Value = "'93.65"
$ws.Range("E17").Value = "  +0.89%  "
$ws.Range("D18").Value = "'0.00001059"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").Value = "'0.06411"
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("E21").Value = "  +1.72%  "
$ws.Range("D22").Value = "'5.795"
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("D23").Value = "'27.815.21"
$ws.Range("E23").Value = "  +1.01%  "
$ws.Range("D24").Value = "'11.30"
$ws.Range("E24").Value = "  +1.41%  "
$ws.Range("D25").Value = "'2.087"
$ws.Range("E25").Value = "  -0.58%  "
$ws.Range("D26").Value = "'164.39"
$ws.Range("E26").Value = "  +1.64%  "
$ws.Range("D27").Value = "'20.28"
$ws.Range("E27").Value = "  -0.75%  "
$ws.Range("D28").Value = "'1.979.85"
$ws.Range("E28").Value = "  +1.33%  "
$ws.Range("E29").Value = "  +4.29%  "
$ws.Range("D30").Value = "'125.92"
$ws.Range("E30").Value = "  +0.89%  "
$ws.Range("E31").Value = "  +0.45%  "
$ws.Range("D32").Value = "'0.09220"
$ws.Range("E32").Value = "  +2.45%  "
$ws.Range("D33").Value = "'3.669"
$ws.Range("E33").Value = "  +0.33%  "
$ws.Range("D34").Value = "'5.532"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").Value = "'11.84"
$ws.Range("E35").Value = "  -1.07%  "
$ws.Range("D36").Value = "'0.02293"
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("D37").Value = "'0.06179"
$ws.Range("E37").Value = "  +2.88%  "
$ws.Range("D38").Value = "'0.2086"
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("D39").Value = "'0.6313"
$ws.Range("E39").Value = "  -0.39%  "
$ws.Range("D40").Value = "'4.945"
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("E41").Value = "  -1.82%  "
$ws.Range("D42").Value = "'1.393"
$ws.Range("E42").Value = "  +0.41%  "
$ws.Range("D43").Value = "'7.786"
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("D44").Value = "'13.27"
$ws.Range("E44").Value = "  +0.89%  "
$ws.Range("D45").Value = "'3.742"
$ws.Range("E45").Value = "  +0.98%  "
$ws.Range("D46").Value = "'0.5888"
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("D47").Value = "'122.26"
$ws.Range("E47").Value = "  +0.23%  "
$ws.Range("D48").Value = "'1.947"
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("D49").Value = "'0.06930"
$ws.Range("E49").Value = "  +1.17%  "
$ws.Range("E50").Value = "  -0.82%  "
$ws.Range("D51").Value = "'72.08"
$ws.Range("E51").Value = "  +0.11%  "
